$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update F2 value from "AY260547.prot" to "angaj2010"
$ws.Range("F2").Value = "angaj2010"

# Update the selection shown in the sheet view to F2
$ws.Activate()
$ws.Range("F2").Select()
